$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target A-column labels (and therefore the per-row container letter used
# in that row comments) are swapped: A<->E, B<->F, C<->G, D<->H.
# Column A cells here are locked, so the sheet must be unprotected for the edit
# and re-protected afterwards to restore the protected state.
$ws.Unprotect()

$ws.Range("A10").Value = "A"
$ws.Range("A11").Value = "B"
$ws.Range("A16").Value = "G"
$ws.Range("A17").Value = "H"
$ws.Range("A22").Value = "C"
$ws.Range("A23").Value = "D"
$ws.Range("A28").Value = "E"
$ws.Range("A29").Value = "F"

$ws.Protect()

# Update the PAML-autogenerated comment text (Container5 -> Container3, and the
# per-cell suffix letter follows the row's new container letter / column index).
$ws.Range("B10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A1")
$ws.Range("C10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A2")
$ws.Range("D10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A3")
$ws.Range("E10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A4")
$ws.Range("F10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A5")
$ws.Range("G10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A6")
$ws.Range("H10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A7")
$ws.Range("I10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A8")
$ws.Range("J10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A9")
$ws.Range("K10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A10")
$ws.Range("L10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A11")
$ws.Range("M10").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_A12")
$ws.Range("B11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B1")
$ws.Range("C11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B2")
$ws.Range("D11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B3")
$ws.Range("E11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B4")
$ws.Range("F11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B5")
$ws.Range("G11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B6")
$ws.Range("H11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B7")
$ws.Range("I11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B8")
$ws.Range("J11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B9")
$ws.Range("K11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B10")
$ws.Range("L11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B11")
$ws.Range("M11").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_B12")
$ws.Range("B16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G1")
$ws.Range("C16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G2")
$ws.Range("D16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G3")
$ws.Range("E16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G4")
$ws.Range("F16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G5")
$ws.Range("G16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G6")
$ws.Range("H16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G7")
$ws.Range("I16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G8")
$ws.Range("J16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G9")
$ws.Range("K16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G10")
$ws.Range("L16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G11")
$ws.Range("M16").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_G12")
$ws.Range("B17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H1")
$ws.Range("C17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H2")
$ws.Range("D17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H3")
$ws.Range("E17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H4")
$ws.Range("F17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H5")
$ws.Range("G17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H6")
$ws.Range("H17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H7")
$ws.Range("I17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H8")
$ws.Range("J17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H9")
$ws.Range("K17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H10")
$ws.Range("L17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H11")
$ws.Range("M17").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_H12")
$ws.Range("B22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C1")
$ws.Range("C22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C2")
$ws.Range("D22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C3")
$ws.Range("E22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C4")
$ws.Range("F22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C5")
$ws.Range("G22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C6")
$ws.Range("H22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C7")
$ws.Range("I22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C8")
$ws.Range("J22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C9")
$ws.Range("K22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C10")
$ws.Range("L22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C11")
$ws.Range("M22").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_C12")
$ws.Range("B23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D1")
$ws.Range("C23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D2")
$ws.Range("D23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D3")
$ws.Range("E23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D4")
$ws.Range("F23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D5")
$ws.Range("G23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D6")
$ws.Range("H23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D7")
$ws.Range("I23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D8")
$ws.Range("J23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D9")
$ws.Range("K23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D10")
$ws.Range("L23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D11")
$ws.Range("M23").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_D12")
$ws.Range("B28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E1")
$ws.Range("C28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E2")
$ws.Range("D28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E3")
$ws.Range("E28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E4")
$ws.Range("F28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E5")
$ws.Range("G28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E6")
$ws.Range("H28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E7")
$ws.Range("I28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E8")
$ws.Range("J28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E9")
$ws.Range("K28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E10")
$ws.Range("L28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E11")
$ws.Range("M28").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_E12")
$ws.Range("B29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F1")
$ws.Range("C29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F2")
$ws.Range("D29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F3")
$ws.Range("E29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F4")
$ws.Range("F29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F5")
$ws.Range("G29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F6")
$ws.Range("H29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F7")
$ws.Range("I29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F8")
$ws.Range("J29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F9")
$ws.Range("K29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F10")
$ws.Range("L29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F11")
$ws.Range("M29").Comment.Text("https://igem.org/Engineering/protocols/Multicolor_particle_calibration/Container3_F12")
